$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Create the new cell style "Normal 2" (font: Aptos Narrow, 11pt, theme color 1)
$normal2 = $wb.Styles.Add("Normal 2")
$normal2.Font.Name = "Aptos Narrow"
$normal2.Font.Size = 11

# Widen column B to fit the longer English prompts
# (target stored width 55.28515625 chars; 54.5 is the input that round-trips
# closest to that value given the engine's column-width quantization)
$ws.Columns.Item(2).ColumnWidth = 54.5

# Populate column C (Brazilian Portuguese translations) for rows 3-70
$ws.Cells.Item(3, 3).Value = 'Escolher uma distribuição Linux para começar usando, ex:  Mint, Ubuntu'
$ws.Cells.Item(3, 3).Style = "Normal 2"
$ws.Cells.Item(4, 3).Value = 'Experimentar o “Windows Subsystem for Linux” (WSL)'
$ws.Cells.Item(4, 3).Style = "Normal 2"
$ws.Cells.Item(5, 3).Value = 'Experimentar o Linux em uma máquina virtual'
$ws.Cells.Item(5, 3).Style = "Normal 2"
$ws.Cells.Item(6, 3).Value = 'Usar o Snappy para baixar aplicativos'
$ws.Cells.Item(6, 3).Style = "Normal 2"
$ws.Cells.Item(7, 3).Value = 'Executar o Linux em um computador antigo'
$ws.Cells.Item(7, 3).Style = "Normal 2"
$ws.Cells.Item(8, 3).Value = 'Baixar ou fazer torrent de uma distribuição Linux'
$ws.Cells.Item(8, 3).Style = "Normal 2"
$ws.Cells.Item(9, 3).Value = 'Encontrar software de código aberto alternativos para usar no Linux'
$ws.Cells.Item(9, 3).Style = "Normal 2"
$ws.Cells.Item(10, 3).Value = 'Instalar o Linux em seu hardware'
$ws.Cells.Item(10, 3).Style = "Normal 2"
$ws.Cells.Item(11, 3).Value = 'Explorar sua nova distribuição Linux'
$ws.Cells.Item(11, 3).Style = "Normal 2"
$ws.Cells.Item(12, 3).Value = 'Fazer dual boot com o Linux em seu hardware'
$ws.Cells.Item(12, 3).Style = "Normal 2"
$ws.Cells.Item(13, 3).Value = 'Usar o Wine para executar aplicativos do Windows no Linux'
$ws.Cells.Item(13, 3).Style = "Normal 2"
$ws.Cells.Item(14, 3).Value = 'Escolher um navegador para usar, ex.: Chromium, Firefox'
$ws.Cells.Item(14, 3).Style = "Normal 2"
$ws.Cells.Item(15, 3).Value = 'Instalar Java ou Python'
$ws.Cells.Item(15, 3).Style = "Normal 2"
$ws.Cells.Item(16, 3).Value = 'Instalar o Steam no Linux'
$ws.Cells.Item(16, 3).Style = "Normal 2"
$ws.Cells.Item(17, 3).Value = 'Instalar e explorar o Fedora Linux'
$ws.Cells.Item(17, 3).Style = "Normal 2"
$ws.Cells.Item(18, 3).Value = 'Usar o Rufus para criar uma imagem de disco'
$ws.Cells.Item(18, 3).Style = "Normal 2"
$ws.Cells.Item(19, 3).Value = 'Configurar uma impressora ou scanner'
$ws.Cells.Item(19, 3).Style = "Normal 2"
$ws.Cells.Item(20, 3).Value = 'Usar o comando man/info para obter ajuda na linha de comando'
$ws.Cells.Item(20, 3).Style = "Normal 2"
$ws.Cells.Item(21, 3).Value = 'Gerenciar as configurações Bluetooth'
$ws.Cells.Item(21, 3).Style = "Normal 2"
$ws.Cells.Item(22, 3).Value = 'Gerenciar as configurações de áudio'
$ws.Cells.Item(22, 3).Style = "Normal 2"
$ws.Cells.Item(23, 3).Value = 'Personalizar sua área de trabalho'
$ws.Cells.Item(23, 3).Style = "Normal 2"
$ws.Cells.Item(24, 3).Value = 'Usar comandos de limpeza do sistema'
$ws.Cells.Item(24, 3).Style = "Normal 2"
$ws.Cells.Item(25, 3).Value = 'Gerenciar diretórios com pwd, cd e ls na linha de comando'
$ws.Cells.Item(25, 3).Style = "Normal 2"
$ws.Cells.Item(26, 3).Value = 'Instalar programas com a linha de comando'
$ws.Cells.Item(26, 3).Style = "Normal 2"
$ws.Cells.Item(27, 3).Value = 'Usar o comando whatis para descobrir o que um comando faz'
$ws.Cells.Item(27, 3).Style = "Normal 2"
$ws.Cells.Item(28, 3).Value = 'Usar o comando kill'
$ws.Cells.Item(28, 3).Style = "Normal 2"
$ws.Cells.Item(29, 3).Value = 'Usar os comandos toilet ou figlet por diversão'
$ws.Cells.Item(29, 3).Style = "Normal 2"
$ws.Cells.Item(30, 3).Value = 'Usar ferramentas básicas que vêm com o Linux, ex: Awk, grep'
$ws.Cells.Item(30, 3).Style = "Normal 2"
$ws.Cells.Item(31, 3).Value = 'Usar um software de backup'
$ws.Cells.Item(31, 3).Style = "Normal 2"
$ws.Cells.Item(32, 3).Value = 'Experimentar um gerenciador de janelas (WM) diferente'
$ws.Cells.Item(32, 3).Style = "Normal 2"
$ws.Cells.Item(33, 3).Value = 'Gerenciar permissões de arquivo com a linha de comando'
$ws.Cells.Item(33, 3).Style = "Normal 2"
$ws.Cells.Item(34, 3).Value = 'Otimizar configurações de gerenciamento de energia'
$ws.Cells.Item(34, 3).Style = "Normal 2"
$ws.Cells.Item(35, 3).Value = 'Usar o comando mount para montar partições'
$ws.Cells.Item(35, 3).Style = "Normal 2"
$ws.Cells.Item(36, 3).Value = 'Usar o Gparted para gerenciar partições de disco'
$ws.Cells.Item(36, 3).Style = "Normal 2"
$ws.Cells.Item(37, 3).Value = 'Configurar o BusyBox em um telefone móvel'
$ws.Cells.Item(37, 3).Style = "Normal 2"
$ws.Cells.Item(38, 3).Value = 'Ler um livro ou assistir tutoriais sobre Linux'
$ws.Cells.Item(38, 3).Style = "Normal 2"
$ws.Cells.Item(39, 3).Value = 'Acidentalmente "brickar" sua máquina com rm -rf'
$ws.Cells.Item(39, 3).Style = "Normal 2"
$ws.Cells.Item(40, 3).Value = 'Resolver um problema por conta própria'
$ws.Cells.Item(40, 3).Style = "Normal 2"
$ws.Cells.Item(41, 3).Value = 'Configurar um Raspberry Pi ou similar'
$ws.Cells.Item(41, 3).Style = "Normal 2"
$ws.Cells.Item(42, 3).Value = 'Usar ssh para acessar remotamente um computador'
$ws.Cells.Item(42, 3).Style = "Normal 2"
$ws.Cells.Item(43, 3).Value = 'Instalar e explorar o Arch Linux'
$ws.Cells.Item(43, 3).Style = "Normal 2"
$ws.Cells.Item(44, 3).Value = 'Instalar e explorar o Kali Linux'
$ws.Cells.Item(44, 3).Style = "Normal 2"
$ws.Cells.Item(45, 3).Value = 'Configurar uma tarefa no cron'
$ws.Cells.Item(45, 3).Style = "Normal 2"
$ws.Cells.Item(46, 3).Value = 'Usar o comando dd'
$ws.Cells.Item(46, 3).Style = "Normal 2"
$ws.Cells.Item(47, 3).Value = 'Instalar um software a partir do código-fonte'
$ws.Cells.Item(47, 3).Style = "Normal 2"
$ws.Cells.Item(48, 3).Value = 'Criar um projeto com um Raspberry Pi ou similar'
$ws.Cells.Item(48, 3).Style = "Normal 2"
$ws.Cells.Item(49, 3).Value = 'Usar Vim e Emacs na linha de comando'
$ws.Cells.Item(49, 3).Style = "Normal 2"
$ws.Cells.Item(50, 3).Value = 'Navegar na internet a partir da linha de comando'
$ws.Cells.Item(50, 3).Style = "Normal 2"
$ws.Cells.Item(51, 3).Value = 'Configurar um servidor de jogos'
$ws.Cells.Item(51, 3).Style = "Normal 2"
$ws.Cells.Item(52, 3).Value = 'Ensinar um amigo a usar o Linux'
$ws.Cells.Item(52, 3).Style = "Normal 2"
$ws.Cells.Item(53, 3).Value = 'Executar Nix ou Guix'
$ws.Cells.Item(53, 3).Style = "Normal 2"
$ws.Cells.Item(54, 3).Value = 'Usar o Linux como seu sistema operacional principal'
$ws.Cells.Item(54, 3).Style = "Normal 2"
$ws.Cells.Item(55, 3).Value = 'Criar ou extrair um arquivo tar usando o terminal'
$ws.Cells.Item(55, 3).Style = "Normal 2"
$ws.Cells.Item(56, 3).Value = 'Instalar o Linux em um Macintosh'
$ws.Cells.Item(56, 3).Style = "Normal 2"
$ws.Cells.Item(57, 3).Value = 'Usar systemctl'
$ws.Cells.Item(57, 3).Style = "Normal 2"
$ws.Cells.Item(58, 3).Value = 'Compilar o kernel Linux a partir do código-fonte'
$ws.Cells.Item(58, 3).Style = "Normal 2"
$ws.Cells.Item(59, 3).Value = 'Instalar o Ubuntu Touch em um dispositivo móvel'
$ws.Cells.Item(59, 3).Style = "Normal 2"
$ws.Cells.Item(60, 3).Value = 'Criar seu próprio serviço de sistema'
$ws.Cells.Item(60, 3).Style = "Normal 2"
$ws.Cells.Item(61, 3).Value = 'Usar o Nmap'
$ws.Cells.Item(61, 3).Style = "Normal 2"
$ws.Cells.Item(62, 3).Value = 'Criar instruções ou um guia de como fazer no Linux'
$ws.Cells.Item(62, 3).Style = "Normal 2"
$ws.Cells.Item(63, 3).Value = 'Usar Tmux para abrir programas e salvar sessões'
$ws.Cells.Item(63, 3).Style = "Normal 2"
$ws.Cells.Item(64, 3).Value = 'Aprender sobre o “filesystem hierarchy standard” (FHS)'
$ws.Cells.Item(64, 3).Style = "Normal 2"
$ws.Cells.Item(65, 3).Value = 'Gerenciar permissões de arquivo estendidas'
$ws.Cells.Item(65, 3).Style = "Normal 2"
$ws.Cells.Item(66, 3).Value = 'Executar o “Linux From Scratch” (LFS)'
$ws.Cells.Item(66, 3).Style = "Normal 2"
$ws.Cells.Item(67, 3).Value = 'Contribuir para a documentação da sua distribuição Linux'
$ws.Cells.Item(67, 3).Style = "Normal 2"
$ws.Cells.Item(68, 3).Value = 'Usar netcat para se comunicar com um servidor'
$ws.Cells.Item(68, 3).Style = "Normal 2"
$ws.Cells.Item(69, 3).Value = 'Dar uma aula sobre Linux'
$ws.Cells.Item(69, 3).Style = "Normal 2"
$ws.Cells.Item(70, 3).Value = 'Executar o Gentoo'
$ws.Cells.Item(70, 3).Style = "Normal 2"

# Update the selection to reflect the newly filled range
$ws.Range("C3:C70").Select()
